$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value = 4274.648500398
$ws.Range("F9").Select()
